$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Teilnahme an Arbeitsgemeinschaften placeholder: ${teilnahme} -> ${ags}
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('${teilnahme}', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$s = $rng.Start
$e = $rng.End

$whole = $d.Range($s, $e)
$whole.Text = '${ags}'

# Force the merged text into three distinct runs (${ / ags / }) by briefly
# toggling formatting on the outer pieces, then reverting it, so the saved
# document keeps the three-run split with matching (identical) run
# properties, mirroring the target markup.
$part1 = $d.Range($s, $s + 2)
$part1.Bold = 1
$part1.Bold = 0

$part3 = $d.Range($s + 5, $s + 6)
$part3.Bold = 1
$part3.Bold = 0

# ---------------------------------------------------------------------------
# 2) Bemerkungen placeholder: ${bemerkungen} -> ${comments_short}
#    (plus a _GoBack bookmark left at the edit point, as Word does)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute('${bemerkungen}', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$s2 = $rng2.Start
$e2 = $rng2.End

$whole2 = $d.Range($s2, $e2)
$whole2.Text = '${comments_short}'

$part1b = $d.Range($s2, $s2 + 2)
$part1b.Bold = 1
$part1b.Bold = 0

$part3b = $d.Range($s2 + 16, $s2 + 17)
$part3b.Bold = 1
$part3b.Bold = 0

$goBackRange = $d.Range($s2 + 16, $s2 + 16)
$d.Bookmarks.Add('_GoBack', $goBackRange)
